$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) cells so numeric-looking strings
# (e.g. "1.001", "10.03") are stored as text, not converted to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '22.324.02'
$ws.Range("E2").Value = '  -4.75%  '
$ws.Range("D3").Value = '1.563.05'
$ws.Range("E3").Value = '  -4.92%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = '289.78'
$ws.Range("E6").Value = '  -3.50%  '
$ws.Range("D7").Value = '0.3734'
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("E8").Value = '  -2.34%  '
$ws.Range("D9").Value = '0.3402'
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("D10").Value = '1.164'
$ws.Range("E10").Value = '  -4.19%  '
$ws.Range("D11").Value = '0.07647'
$ws.Range("E11").Value = '  -5.11%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '21.42'
$ws.Range("E13").Value = '  -3.01%  '
$ws.Range("D14").Value = '6.026'
$ws.Range("E14").Value = '  -4.20%  '
$ws.Range("D15").Value = '6.915'
$ws.Range("E15").Value = '  -4.60%  '
$ws.Range("D16").Value = '1.559.53'
$ws.Range("E16").Value = '  -5.32%  '
$ws.Range("D17").Value = '0.00001127'
$ws.Range("E17").Value = '  -6.87%  '
$ws.Range("D18").Value = '89.80'
$ws.Range("E18").Value = '  -5.76%  '
$ws.Range("D19").Value = '0.06719'
$ws.Range("E19").Value = '  -3.91%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '6.235'
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").Value = '16.56'
$ws.Range("E22").Value = '  -4.76%  '
$ws.Range("D23").Value = '0.5281'
$ws.Range("E23").Value = '  -7.99%  '
$ws.Range("E24").Value = '  -3.69%  '
$ws.Range("D25").Value = '22.320.50'
$ws.Range("D26").Value = '2.401'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").Value = '2.806'
$ws.Range("E27").Value = '  -6.67%  '
$ws.Range("D28").Value = '20.17'
$ws.Range("E28").Value = '  -3.98%  '
$ws.Range("D29").Value = '146.01'
$ws.Range("E29").Value = '  -3.84%  '
$ws.Range("D30").Value = '4.982'
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("D31").Value = '125.27'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").Value = '1.734.81'
$ws.Range("E32").Value = '  -5.42%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.165'
$ws.Range("E33").Value = '  -10.01%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '1.007'
$ws.Range("E34").Value = '  +1.99%  '
$ws.Range("D35").Value = '2.013'
$ws.Range("E35").Value = '  -6.16%  '
$ws.Range("D36").Value = '10.03'
$ws.Range("E36").Value = '  -10.11%  '
$ws.Range("D37").Value = '0.08500'
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("D38").Value = '0.02551'
$ws.Range("E38").Value = '  -5.12%  '
$ws.Range("D39").Value = '0.2312'
$ws.Range("E39").Value = '  -4.31%  '
$ws.Range("E40").Value = '  -7.18%  '
$ws.Range("D41").Value = '1.310'
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("D42").Value = '0.06388'
$ws.Range("E42").Value = '  -5.74%  '
$ws.Range("D43").Value = '11.69'
$ws.Range("E43").Value = '  -9.08%  '
$ws.Range("D44").Value = '0.6341'
$ws.Range("E44").Value = '  -7.82%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '14.01'
$ws.Range("E46").Value = '  -9.85%  '
$ws.Range("D47").Value = '0.5966'
$ws.Range("E47").Value = '  -6.55%  '
$ws.Range("E48").Value = '  -4.36%  '
$ws.Range("D49").Value = '2.086'
$ws.Range("E49").Value = '  -6.79%  '
$ws.Range("D50").Value = '1.264'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").Value = '124.56'
$ws.Range("E51").Value = '  -2.25%  '

# Reset number format back to General/Normal style so cells retain their
# original (unstyled) appearance, matching the source workbook.
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}